# Trade #39 closed at 2026-02-16 22:56:00 - base_strategy UP +0.000%
#
# Append a new trade row (row 40) to both the "All Trades" sheet and the
# "base_strategy" sheet with identical data. The new trade shares its
# Date / Entry Price / Status / Entry Reason with the previous row (39),
# so that row is duplicated via Copy (preserving cell types/formatting
# exactly) and only the cells that actually differ are overwritten.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "base_strategy")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Duplicate the last existing trade row (39) into the new row (40) -
    # this keeps formatting/types (e.g. the text-formatted Date column)
    # identical to the rest of the table.
    $ws.Range("A39:Q39").Copy($ws.Range("A40:Q40"))

    # Overwrite only the fields that differ for trade #39.
    $ws.Cells.Item(40, 1).Value = 39              # Trade #
    $ws.Cells.Item(40, 3).Value = "22:56:00"      # Time
    $ws.Cells.Item(40, 5).Value = "UP"            # Side
}
